$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# --- Title shape ("Title 1") ---
$title = $s.Shapes.Item(1)
$title.Left = 28.956536333070876
$title.Width = 655.0434875669213
$title.Height = 54.02622047244095
$title.TextFrame.TextRange.Text = "PM Probes in Loopback Mode for SR Policy"

# --- Content Placeholder 2 (body bullets) ---
$content = $s.Shapes.Item(3)
$content.Height = 151.31614173228346

# --- Rectangle 2 (diagram box with t1/t4 figure) ---
$rect = $s.Shapes.Item(4)
$rect.Left = 211.5270920141729

# Update the figure caption run in place (keep it a single run)
$caption = $rect.TextFrame.TextRange.Paragraphs(13).Runs(1)
$caption.Text = "           Figure: PM Loopback Mode"
